# excel-app/DataBase.xlsx  -  "Rename Fetch to DataBase - Add worksheet generation"
#
# 1) Well sheet: insert a new column B "UWI" populated with the well's
#    UWI string (quote-prefixed, like a typed '100 14 36 085 17 W3 00).
# 2) Monthly sheet: insert two new columns before the old RTPPrice column,
#    label them SalesPrice / TransPrice, and fill the data rows with a
#    constant SalesPrice (100) and TransPrice (2.2) pair.
# 3) Restore/point the window's first-visible-sheet back near the start,
#    and leave the final on-screen selections matching the authored state
#    (Well!E1 and Monthly!L3:L37), with Monthly as the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Well" worksheet - insert UWI column (new column B)
# ---------------------------------------------------------------------
$wsWell = $wb.Worksheets.Item("Well")

$wsWell.Columns("B:B").Insert()
$wsWell.Range("B1").Value = "UWI"

$uwi = "100 14 36 085 17 W3 00"
$lastRow = 26
for ($r = 2; $r -le $lastRow; $r++) {
    # Leading apostrophe = "entered as text" (quote-prefixed), matching
    # the workbook's existing style used for similar look-alike-numeric
    # text values.
    $wsWell.Cells.Item($r, 2).Value = "'" + $uwi
}

# ---------------------------------------------------------------------
# 2. "Monthly" worksheet - insert SalesPrice / TransPrice columns
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly")

$wsMonthly.Columns("K:L").Insert()
$wsMonthly.Range("K1").Value = "SalesPrice"
$wsMonthly.Range("L1").Value = "TransPrice"

$lastDataRow = 37
$wsMonthly.Range("K2:K" + $lastDataRow).Value = 100
$wsMonthly.Range("L2:L" + $lastDataRow).Value = 2.2

# ---------------------------------------------------------------------
# 3. View state - first visible sheet / selections / active sheet
# ---------------------------------------------------------------------
# Scroll the tab strip back toward the front (best-effort; matches the
# authored workbookView firstSheet move from index 3 down to index 1).
$excel.ActiveWindow.ScrollWorkbookTabs(1, 1)

# Leave Well's selection on the new header cell...
$wsWell.Range("E1").Select()

# ...and finish on Monthly (making it the active tab again) with its
# selection over the freshly filled TransPrice column.
$wsMonthly.Range("L3:L37").Select()
